$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: account holder name and card number
$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit card number that must remain stored as text (not a
# number) while keeping its original cell style. Writing the value directly
# (even with a leading apostrophe) causes Excel to stamp a new "quote
# prefix" style on the cell. Instead, stage the text value in a scratch
# cell formatted as Text, copy it, and paste only the value into B3 so the
# destination keeps its original style/format.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "2570314725427075"
$ws.Range("Z1").Copy()
$ws.Range("B3").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").Clear()

$ws.Range("C3").Value = "Mohaupt"

# Row 5: starting balance date
$ws.Range("D5").Value = "KONTOSTAND AM 25.03.2024"

# Row 6
$ws.Range("B6").Value = "26.03."
$ws.Range("C6").Value = "27.03."
$ws.Range("D6").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E6").Value = "24,87-"

# Row 7
$ws.Range("B7").Value = "28.03."
$ws.Range("C7").Value = "29.03."
$ws.Range("D7").Value = "AMAZON.DE MKTPLC EU RNAAEA"
$ws.Range("E7").Value = "233,90-"

# Row 8
$ws.Range("B8").Value = "29.03."
$ws.Range("C8").Value = "30.03."
$ws.Range("D8").Value = "BEITRAG Allianz SE K-57342728"
$ws.Range("E8").Value = "57,36-"

# Row 9
$ws.Range("B9").Value = "02.04."
$ws.Range("C9").Value = "03.04."
$ws.Range("D9").Value = "RECHNUNG VODAFONE GMBH 43123452"
$ws.Range("E9").Value = "41,95-"

# Row 10
$ws.Range("B10").Value = "05.04."
$ws.Range("C10").Value = "06.04."
$ws.Range("D10").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 40290511"
$ws.Range("E10").Value = "87,67-"

# Row 11
$ws.Range("B11").Value = "08.04."
$ws.Range("C11").Value = "09.04."
$ws.Range("D11").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E11").Value = "50,06-"

# Row 12: ending balance date and amount
$ws.Range("D12").Value = "KONTOSTAND AM 12.04.2024"
$ws.Range("E12").Value = "495,81-"

# Row 13: next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 21.04.2024"
